# Fruta / hortaliza, semanal
# Inserts a new weekly price observation as row 41 (pushing the existing
# rows 41-44 down to 42-45) on the "Vega Modelo de Temuco - Locoto" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 41; everything below shifts
# down one row (old 41 -> 42, old 42 -> 43, old 43 -> 44, old 44 -> 45).
$ws.Rows("41:41").Insert()

# Populate the newly inserted row 41 with the new observation.
$ws.Range("A41").Value = 10
$ws.Range("B41").Value = "Vega Modelo de Temuco"
$ws.Range("C41").Value = "La Araucanía"
$ws.Range("D41").Value = 44669
$ws.Range("E41").Value = 9
$ws.Range("F41").Value = 100112042
$ws.Range("G41").Value = "Locoto"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 60
$ws.Range("K41").Value = 6250
$ws.Range("L41").Value = 6250
$ws.Range("M41").Value = 6250
$ws.Range("N41").Value = "$/kilo"
$ws.Range("O41").Value = "Región de Arica y Parinacota"
$ws.Range("P41").Value = 6250
$ws.Range("Q41").Value = 1
$ws.Range("R41").Value = "Hortaliza"
